# First version of Overall CI Arch
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# 1) "Artifact Store" -> "Object Store"
$s.Shapes.Item(3).TextFrame.TextRange.Text = "Object Store"

# 2) Resize/reposition "Container Image Registry" flowchart shape
$shRegistry = $s.Shapes.Item(11)
$shRegistry.Left = 760.9314960629921
$shRegistry.Top = 266.5127559055118
$shRegistry.Width = 102.73701187401575
$shRegistry.Height = 30.643938007874016

# 3) "Docker" -> "Docker " + "Cli" (two runs)
$shDocker = $s.Shapes.Item(45)
$trDocker = $shDocker.TextFrame.TextRange
$trDocker.Text = "Docker "
$null = $trDocker.InsertAfter("Cli")

# 4) Shrink the down arrow under Docker
$shDownArrow = $s.Shapes.Item(46)
$shDownArrow.Height = 27.33488188976378

# 5) Reposition/resize the bent-up arrow
$shBentArrow = $s.Shapes.Item(57)
$shBentArrow.Left = 827.1697697795275
$shBentArrow.Top = 281.6078042755905
$shBentArrow.Width = 24.384566929133857
$shBentArrow.Height = 62.95370078740157

# 6) Add a new "Docker hub" rounded-rectangle shape after "CNI (Calico)",
#    by duplicating it so style/line/fill refs + body formatting carry over.
$shCni = $s.Shapes.Item(60)
$dupRange = $shCni.Duplicate()
$newShape = $dupRange.Item(1)
$newShape.Name = "Rectangle: Rounded Corners 78"
$newShape.Left = 760.9314960629921
$newShape.Top = 250.82291438582678
$newShape.Width = 102.12023622047244
$newShape.Height = 14.846141732283465
$newShape.Fill.ForeColor.RGB = 5287936
$newShape.TextFrame.TextRange.Text = "Docker hub"
$newShape.TextFrame.TextRange.Font.Size = 12
